$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in B1 and C1
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$ws.Range("B1").Value2 = $c1
$ws.Range("C1").Value2 = $b1

# Swap B/C values for rows 2 through 16 (the two data columns were swapped)
for ($r = 2; $r -le 16; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $bVal = $bCell.Value2
    $cVal = $cCell.Value2
    $bCell.Value2 = $cVal
    $cCell.Value2 = $bVal
}

# Row 16's date cell currently carries the special "last row" number format
# (YYYY-MM-DD). Row 17 becomes the new last row, so it should inherit that
# format, while row 16 switches to the regular date format used by the
# preceding rows (YYYY-MM-DD HH:MM:SS).
$lastRowFormat = $ws.Cells.Item(16, 1).NumberFormat
$ws.Cells.Item(16, 1).NumberFormat = $ws.Cells.Item(15, 1).NumberFormat

# Add new row 17 with the new final data point, using the number format
# previously carried by row 16.
$ws.Cells.Item(17, 1).Value2 = 45749
$ws.Cells.Item(17, 1).NumberFormat = $lastRowFormat
$ws.Cells.Item(17, 2).Value2 = 764.244
$ws.Cells.Item(17, 3).Value2 = 800.9299999999999
